# Apply a cyclic rotation of the data in rows 2-4 for columns A,B,E,F,G,H,Q,R
# after.row2 = before.row3
# after.row3 = before.row4
# after.row4 = before.row2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

# Capture the original values for rows 2, 3, 4 before overwriting anything.
$orig = @{}
foreach ($row in 2..4) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value2
    }
}

# Cyclic shift: row2 <- row3, row3 <- row4, row4 <- row2
$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($row in 2..4) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $orig[$srcRow][$col]
    }
}
